$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.015751481056213
$ws.Range("B1").Value = 1.341597676277161
$ws.Range("C1").Value = 2.145710945129395
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.996834993362427
